$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from row 61 (A61 -> style 13 date border7, B61 -> style 14 border8)
# into the currently-empty row 62, then set values.
$ws.Range("A61:B61").Copy()
$ws.Range("A62:B62").PasteSpecial(-4122)

$ws.Range("A62").Value = 45179
$ws.Range("B62").Value = ""

# Row 63 already has D63; copy formats from row61 A/B into A63/B63 too
$ws.Range("A63:B63").PasteSpecial(-4122)
$ws.Range("A63").Value = 45180
$ws.Range("B63").Value = "fix TfidfProcessor Preprocessing (punctations into ' ')"

$ws.Rows.Item(62).RowHeight = 17
$ws.Rows.Item(63).RowHeight = 18

$ws.Range("A64").Select()

Write-Host "done"
